$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mousebites removed / gerbers + drill files re-exported; fiducial INSERT set
# was missing before, so export the XY position/rotation data for the board
# fiducials (3 on top, 3 on bottom) as new rows appended to the pick & place
# style table.
$fiducials = @(
    @("FD1", 2.86, 2.0099999999999998, "Top",    0, "FIDUCIAL"),
    @("FD2", 0,    0,                  "Top",    0, "FIDUCIAL"),
    @("FD3", 2.86, 0,                  "Top",    0, "FIDUCIAL"),
    @("FD4", 0,    0,                  "Bottom", 0, "FIDUCIAL"),
    @("FD5", 2.84, 0,                  "Bottom", 0, "FIDUCIAL"),
    @("FD6", 2.84, 2.0099999999999998, "Bottom", 0, "FIDUCIAL")
)

$startRow = 99
for ($i = 0; $i -lt $fiducials.Count; $i++) {
    $row = $startRow + $i
    $data = $fiducials[$i]
    $ws.Range("A$row").Value = $data[0]
    $ws.Range("B$row").Value = $data[1]
    $ws.Range("C$row").Value = $data[2]
    $ws.Range("D$row").Value = $data[3]
    $ws.Range("E$row").Value = $data[4]
    $ws.Range("F$row").Value = $data[5]
}

# Widen the MPN/Notes column (F) now that it also carries "FIDUCIAL" labels
$ws.Columns.Item(6).ColumnWidth = 8.14

# Move the selection/scroll position down to the newly appended data
[void]$ws.Range("D$($startRow + $fiducials.Count - 1)").Select()
